# Update the cosinor analysis results for rows 2 and 3 with re-run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"22.78000000000012"
$ws.Range("G2").Value = [double]"1.319934592181049e-06"
$ws.Range("H2").Value = [double]"2.838834339187445e-05"
$ws.Range("K2").Value = [double]"5.200072658106935"
$ws.Range("L2").Value = "[2.6910732551797727, 7.709072061034097]"
$ws.Range("M2").Value = [double]"5.599489521146417e-05"
$ws.Range("N2").Value = [double]"0.0001119897904229283"
$ws.Range("P2").Value = "[-1.6855792415656943, -0.6289474781961548]"
$ws.Range("Q2").Value = [double]"2.108583131210651e-05"
$ws.Range("R2").Value = [double]"4.217166262421301e-05"
$ws.Range("S2").Value = [double]"11.47659430448781"
$ws.Range("T2").Value = "[10.107992554409417, 12.845196054566209]"
$ws.Range("W2").Value = [double]"4.195715715715739"
$ws.Range("X2").Value = [double]"2.280280280280294"
$ws.Range("Y2").Value = [double]"6.111151151151184"

# Row 3 updates
$ws.Range("E3").Value = [double]"24.22000000000035"
$ws.Range("G3").Value = [double]"0.0005084099697493238"
$ws.Range("H3").Value = [double]"0.001769037968785438"
$ws.Range("K3").Value = [double]"5.03220610686169"
$ws.Range("L3").Value = "[2.054394084207873, 8.010018129515506]"
$ws.Range("M3").Value = [double]"0.0009955938694530264"
$ws.Range("N3").Value = [double]"0.0009955938694530264"
$ws.Range("P3").Value = "[0.723289599925578, 2.23276354759635]"
$ws.Range("Q3").Value = [double]"0.0001432777040062394"
$ws.Range("R3").Value = [double]"0.0001432777040062394"
$ws.Range("S3").Value = [double]"12.19574729301942"
$ws.Range("T3").Value = "[10.423711658361967, 13.967782927676872]"
$ws.Range("W3").Value = [double]"18.52260260260287"
$ws.Range("X3").Value = [double]"15.61329329329351"
$ws.Range("Y3").Value = [double]"21.43191191191222"

$wb.Save()
Write-Output "updated cosinor results for rows 2 and 3"
